$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "59.352.03"
$ws.Range("E2").Value = "  +0.89%  "

# Row 3
$ws.Range("D3").Value = "2.590.48"
$ws.Range("E3").Value = "  -0.26%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.72%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.56%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.03%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.603"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.16%  "

# Row 9
$ws.Range("D9").Value = "2.599.82"
$ws.Range("E9").Value = "  -0.10%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.67"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.53%  "

# Row 11
$ws.Range("E11").Value = "  +3.43%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.158"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +11.36%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.347"
$ws.Range("D13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "3.042.89"
$ws.Range("E14").Value = "  -0.42%  "

# Row 15
$ws.Range("D15").Value = "59.370.95"
$ws.Range("E15").Value = "  +1.00%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.58"
$ws.Range("D16").Style = "Normal"

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000138"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.31%  "

# Row 18
$ws.Range("D18").Value = "2.597.14"
$ws.Range("E18").Value = "  -0.17%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.61%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "338.50"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.33%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.90%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.21"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.63%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.21%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.54"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.45%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.455"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.45%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.49%  "

# Row 27
$ws.Range("E27").Value = "  +1.92%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.29"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.32%  "

# Row 29
$ws.Range("D29").Value = "0.0₃0783"
$ws.Range("E29").Value = "  +3.97%  "

# Row 30
$ws.Range("E30").Value = "  +0.00%  "

# Row 31
$ws.Range("E31").Value = "  +1.00%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.09"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.66%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "159.25"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.02%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.09"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.74%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.04"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.97%  "

# Row 36
$ws.Range("E36").Value = "  +1.71%  "

# Row 37
$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.884"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.38%  "

# Row 38
$ws.Range("B38").Value = "SuiNetwork"
$ws.Range("C38").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.879"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.67%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.11"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.57%  "

# Row 40
$ws.Range("E40").Value = "  +2.66%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "295.32"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.01%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.68"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.18%  "

# Row 43
$ws.Range("E43").Value = "  +0.09%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0979"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.47%  "

# Row 45
$ws.Range("E45").Value = "  -0.27%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0539"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.39%  "

# Row 47
$ws.Range("E47").Value = "  +3.09%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.62"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.10%  "

# Row 49
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0233"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.55%  "

# Row 50
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "124.39"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.06%  "

# Row 51
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.56"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.65%  "

